$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '29.991.65'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '1.878.63'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('D4').Value = "'" + '0.9996'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'" + '242.85'
$ws.Range('E5').Value = '  -3.43%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = "'" + '0.4953'
$ws.Range('D8').Value = "'" + '0.2919'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('E9').Value = '  -2.07%  '
$ws.Range('D10').Value = '1.882.32'
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').Value = "'" + '16.77'
$ws.Range('E11').Value = '  -2.67%  '
$ws.Range('D12').Value = "'" + '0.07249'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = "'" + '0.6674'
$ws.Range('E13').Value = '  -3.27%  '
$ws.Range('D14').Value = "'" + '86.64'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').Value = "'" + '4.910'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = '29.966.03'
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').Value = "'" + '0.000007875'
$ws.Range('E17').Value = '  -3.25%  '
$ws.Range('D18').Value = "'" + '0.9994'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = "'" + '12.78'
$ws.Range('E19').Value = '  -1.43%  '
$ws.Range('D20').Value = '2.122.95'
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').Value = "'" + '0.9996'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = "'" + '4.776'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('D23').Value = "'" + '5.744'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').Value = "'" + '9.074'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'" + '149.72'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('B26').Value = 'BitcoinCash'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D26').Value = "'" + '142.12'
$ws.Range('E26').Value = '  +4.81%  '
$ws.Range('D27').Value = "'" + '17.03'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = "'" + '1.921'
$ws.Range('E28').Value = '  -3.57%  '
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').Value = "'" + '4.196'
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('D31').Value = "'" + '0.08751'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').Value = "'" + '3.971'
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('D33').Value = "'" + '0.05069'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').Value = "'" + '1.118'
$ws.Range('E34').Value = '  -2.35%  '
$ws.Range('D35').Value = "'" + '0.7125'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = "'" + '0.01803'
$ws.Range('E37').Value = '  +6.55%  '
$ws.Range('D38').Value = "'" + '2.690'
$ws.Range('E38').Value = '  -4.28%  '
$ws.Range('D39').Value = "'" + '2.176'
$ws.Range('E39').Value = '  -4.02%  '
$ws.Range('D40').Value = "'" + '0.9324'
$ws.Range('E40').Value = '  -3.74%  '
$ws.Range('D41').Value = "'" + '5.779'
$ws.Range('E41').Value = '  -5.84%  '
$ws.Range('D42').Value = "'" + '0.4242'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D44').Value = "'" + '102.66'
$ws.Range('E44').Value = '  -2.18%  '
$ws.Range('D45').Value = "'" + '7.448'
$ws.Range('E45').Value = '  -1.94%  '
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('D48').Value = "'" + '32.59'
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('D49').Value = "'" + '8.326'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').Value = "'" + '0.3787'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').Value = "'" + '55.97'
$ws.Range('E51').Value = '  -1.25%  '
